$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "구글 코랩(Google Colab)에서 Mecab 형태소 분석기, konlpy 쉽게 설치하기"
$ws.Range("E4").Value = "https://teddylee777.github.io/colab/colab-mecab"

$ws.Range("D9").Value = "자율주행차에 게임이론이 들어간다고?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/auto-driving-game-theory/#utm_source=rss&utm_medium=rss&utm_campaign=auto-driving-game-theory"

$ws.Range("D10").Value = "블록체인 노트"
$ws.Range("E10").Value = "https://rokrokss.com/post/2022/01/13/%EB%B8%94%EB%A1%9D%EC%B2%B4%EC%9D%B8-%EB%85%B8%ED%8A%B8.html"

$ws.Range("D28").Value = "[Manipulator] Control"
$ws.Range("E28").Value = "https://ropiens.tistory.com/176"

$ws.Range("D46").Value = "코로나예방접종 후, 심근염 및 심낭염"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/426"

$ws.Range("D51").Value = "[윈도우10] 메모장을 관리자 권한으로 열어야 할 때"
$ws.Range("E51").Value = "https://bskyvision.com/1238"
